$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed each new row from row 14s formatting (font/numberformat), then overwrite values
# and add hyperlinks on column C, matching the source workbook pattern (rows 2-14).

$ws.Range('A14:F14').Copy($ws.Range('A15'))
$ws.Range('A15').Value = 'A047'
$ws.Range('B15').Value = '경기도_김포시'
$ws.Range('C15').Value = 'https://www.gimpo.go.kr/portal/ntfcPblancList.do?key=1004&cate_cd=1&searchCnd=40900000000&pageUnit=90'
$ws.Hyperlinks.Add($ws.Range('C15'), 'https://www.gimpo.go.kr/portal/ntfcPblancList.do?key=1004&cate_cd=1&searchCnd=40900000000&pageUnit=90') | Out-Null
$ws.Range('D15').Value = '2024년 김포시 미술작품 구입 심의위원회 결과 공고'
$ws.Range('E15').Value = 45656.0
$ws.Range('F15').Value = 45656.839004629626
$ws.Range('E15:F15').NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Range('A14:F14').Copy($ws.Range('A16'))
$ws.Range('A16').Value = 'A047'
$ws.Range('B16').Value = '경기도_김포시'
$ws.Range('C16').Value = 'https://www.gimpo.go.kr/portal/ntfcPblancList.do?key=1004&cate_cd=1&searchCnd=40900000000&pageUnit=90'
$ws.Hyperlinks.Add($ws.Range('C16'), 'https://www.gimpo.go.kr/portal/ntfcPblancList.do?key=1004&cate_cd=1&searchCnd=40900000000&pageUnit=90') | Out-Null
$ws.Range('D16').Value = '2025년도 1분기 도로관리심의회 사업계획서 제출 알림 공고'
$ws.Range('E16').Value = 45656.0
$ws.Range('F16').Value = 45656.839004629626
$ws.Range('E16:F16').NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Range('A14:F14').Copy($ws.Range('A17'))
$ws.Range('A17').Value = 'A059'
$ws.Range('B17').Value = '경기도_안성시'
$ws.Range('C17').Value = 'https://www.anseong.go.kr/portal/saeol/gosiList.do?mId=0501040000&token=1717572030185'
$ws.Hyperlinks.Add($ws.Range('C17'), 'https://www.anseong.go.kr/portal/saeol/gosiList.do?mId=0501040000&token=1717572030185') | Out-Null
$ws.Range('D17').Value = '안성시 투자유치심의위원회 위원 공개모집 공고'
$ws.Range('E17').Value = 45656.0
$ws.Range('F17').Value = 45656.839004629626
$ws.Range('E17:F17').NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Range('A14:F14').Copy($ws.Range('A18'))
$ws.Range('A18').Value = 'A126'
$ws.Range('B18').Value = '전라도_전주시'
$ws.Range('C18').Value = 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A'
$ws.Hyperlinks.Add($ws.Range('C18'), 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A') | Out-Null
$ws.Range('D18').Value = '제안서 평가결과 공개(2025 시정 소식지 전주다움 제작)'
$ws.Range('E18').Value = 45656.0
$ws.Range('F18').Value = 45656.839004629626
$ws.Range('E18:F18').NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Range('A14:F14').Copy($ws.Range('A19'))
$ws.Range('A19').Value = 'A126'
$ws.Range('B19').Value = '전라도_전주시'
$ws.Range('C19').Value = 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A'
$ws.Hyperlinks.Add($ws.Range('C19'), 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A') | Out-Null
$ws.Range('D19').Value = '제안서 평가결과 공개(2025 시정 소식지 전주다움 제작)'
$ws.Range('E19').Value = 45656.0
$ws.Range('F19').Value = 45656.839004629626
$ws.Range('E19:F19').NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Range('A14:F14').Copy($ws.Range('A20'))
$ws.Range('A20').Value = 'A175'
$ws.Range('B20').Value = '경상도_울진군'
$ws.Range('C20').Value = 'https://www.uljin.go.kr/index.uljin?menuCd=DOM_000000103002007000'
$ws.Hyperlinks.Add($ws.Range('C20'), 'https://www.uljin.go.kr/index.uljin?menuCd=DOM_000000103002007000') | Out-Null
$ws.Range('D20').Value = '직산1리 연안재해방지시설 설치사업”특정공법(특허・신기술) 선정을 위한 공법선정위원회 결과 공개'
$ws.Range('E20').Value = 45656.0
$ws.Range('F20').Value = 45656.839004629626
$ws.Range('E20:F20').NumberFormat = 'yyyy-mm-dd h:mm:ss'
